# Apply the ValueSet-source-classification-values.xlsx metadata update:
#  - Rename "Include from Primary or Seconda" sheet to "Include from Classification o"
#  - Update Title, Date, Description on the Metadata sheet
#  - Update the System URI on the "Include from ..." sheet

$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Title (A5/B5): "Primary or Secondary Value Set" -> "Source Classification Value Set"
$meta.Range("B5").Value = "Source Classification Value Set"

# Date (A8/B8): refreshed timestamp
$meta.Range("B8").Value = "2021-10-01T15:07:10+00:00"

# Description (A12/B12): "Source Classification Value Set" -> "Value set for classifying data origin"
$meta.Range("B12").Value = "Value set for classifying data origin"

# --- Include from ... sheet ---
$includeSheet = $wb.Worksheets.Item("Include from Primary or Seconda")

# System URI (A4/B4)
$includeSheet.Range("B4").Value = "http://ibm.com/fhir/cdm/CodeSystem/process-meta-source-classification"

# Rename the sheet itself last so the lookup above still works
$includeSheet.Name = "Include from Classification o"
